# Update column G ("K") values on Sheet1 with freshly regenerated strikeout counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> new K value
$values = @{
    2  = 4
    3  = 2
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 3
    17 = 2
    18 = 2
    19 = 1
    20 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
